# Revert "FInishing for first demo":
# remove the verb data that was filled in (A8, and the big B123:F166 block
# on Sheet1), and normalize the couple of cells that carried the now-unused
# duplicate "Arial/family=2" font (styles.xml fontId 4 / cellXfs 4).

$wb = $excel.ActiveWorkbook

# ---- Sheet1 : clear the reverted verb rows -------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# row 8 lost its "doe" stam entry
$ws1.Range("A8").ClearContents()

# rows 123-166 (staan..zwijgen) lost all their Infinitief/Impefectum/Zijn/
# Pefectum/English content
$ws1.Range("A123:F166").ClearContents()

# C165/D165 used to carry the stray duplicate-font style (cellXfs index 4);
# after the data is gone, bring them back in line with the rest of the
# column (copy the plain formatting from a neighbouring normal cell).
$ws1.Range("E165").Copy()
$ws1.Range("C165:D165").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Sheet1's column C widened back out
$ws1.Columns("C").ColumnWidth = 23.833333333333336

# ---- Sheet7 : drop the stray duplicate-font style on B27:C27 -------------
$ws7 = $wb.Worksheets.Item("Sheet7")
$ws7.Range("D27").Copy()
$ws7.Range("B27:C27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Cosmetic column-width tweaks on the per-letter sheets ---------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Columns("B").ColumnWidth = 14.966666666666667

$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Columns("B").ColumnWidth = 24.333333333333336

# ---- Restore the active selection on Sheet1 ------------------------------
$ws1.Activate()
$ws1.Range("B123:C123").Select()
